$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "`$RAND('-1')example@example.com"
$ws.Range("C5").Value = "`$RAND('0')example@example.com"
$ws.Range("C6").Value = "`$RAND('1')example@example.com"
$ws.Range("C7").Value = "`$RAND('2')example@example.com"
$ws.Range("C8").Value = "`$RAND('10')example@example.com"

$ws.Range("C9").Select()
